$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to match new data source (real data column headers)
$ws.Range("A1").Value = "angular frequency"
$ws.Range("B1").Value = "Z_imag"
$ws.Range("C1").Value = "Z_real"
$ws.Range("D1").Value = "applied voltage"
$ws.Range("E1").Value = "J_ph"
$ws.Range("F1").Value = "J"
$ws.Range("G1").Value = "abs(Z)"
$ws.Range("H1").Value = "theta"

# Update selection to reflect new active cell
$ws.Range("C1").Select()
